# Update countries & provincias Spain
#
# The source diff updates the COVID-19 country table with newer figures
# for "10 de Mayo de 2020". Most rows simply get refreshed numbers, but
# two pairs of countries also swap their rank/position in the shared
# string table because their "Casos totales" (column B) crossed over:
#   - Chile (row 24) overtook Ecuador (row 25)
#   - "Consejo Danes para los Refugiados" (row 97) overtook
#     Mayotte (row 98) / Kirguistan (row 99)
# Every cell below is written with its final, post-edit value (label +
# the 7 numeric columns B:H) so both the label reshuffle and the data
# refresh happen together per-row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 1370999
$ws.Cells.Item(4, 3).Value = 3361
$ws.Cells.Item(4, 4).Value = 256910
$ws.Cells.Item(4, 5).Value = 1033219
$ws.Cells.Item(4, 6).Value = 16514
$ws.Cells.Item(4, 7).Value = 83
$ws.Cells.Item(4, 8).Value = 80870

# Row 24: Chile (was Ecuador; Chile now ranks here)
$ws.Cells.Item(24, 1).Value = "Chile"
$ws.Cells.Item(24, 2).Value = 30063
$ws.Cells.Item(24, 3).Value = 1197
$ws.Cells.Item(24, 4).Value = 13605
$ws.Cells.Item(24, 5).Value = 16135
$ws.Cells.Item(24, 6).Value = 574
$ws.Cells.Item(24, 7).Value = 11
$ws.Cells.Item(24, 8).Value = 323

# Row 25: Ecuador (was Chile; keeps the old Ecuador figures)
$ws.Cells.Item(25, 1).Value = "Ecuador"
$ws.Cells.Item(25, 2).Value = 29559
$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(25, 4).Value = 3433
$ws.Cells.Item(25, 5).Value = 23999
$ws.Cells.Item(25, 6).Value = 181
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 2127

# Row 29: Singapur
$ws.Cells.Item(29, 1).Value = "Singapur"
$ws.Cells.Item(29, 2).Value = 23822
$ws.Cells.Item(29, 3).Value = 486
$ws.Cells.Item(29, 4).Value = 3225
$ws.Cells.Item(29, 5).Value = 20576
$ws.Cells.Item(29, 6).Value = 24
$ws.Cells.Item(29, 7).Value = 1
$ws.Cells.Item(29, 8).Value = 21

# Row 52: Noruega
$ws.Cells.Item(52, 1).Value = "Noruega"
$ws.Cells.Item(52, 2).Value = 8122
$ws.Cells.Item(52, 3).Value = 17
$ws.Cells.Item(52, 4).Value = 32
$ws.Cells.Item(52, 5).Value = 7871
$ws.Cells.Item(52, 6).Value = 22
$ws.Cells.Item(52, 7).Value = 0
$ws.Cells.Item(52, 8).Value = 219

# Row 58: Argelia
$ws.Cells.Item(58, 1).Value = "Argelia"
$ws.Cells.Item(58, 2).Value = 5891
$ws.Cells.Item(58, 3).Value = 168
$ws.Cells.Item(58, 4).Value = 2841
$ws.Cells.Item(58, 5).Value = 2543
$ws.Cells.Item(58, 6).Value = 22
$ws.Cells.Item(58, 7).Value = 5
$ws.Cells.Item(58, 8).Value = 507

# Row 70: Irak
$ws.Cells.Item(70, 1).Value = "Irak"
$ws.Cells.Item(70, 2).Value = 2818
$ws.Cells.Item(70, 3).Value = 51
$ws.Cells.Item(70, 4).Value = 1790
$ws.Cells.Item(70, 5).Value = 918
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 1
$ws.Cells.Item(70, 8).Value = 110

# Row 71: Grecia
$ws.Cells.Item(71, 1).Value = "Grecia"
$ws.Cells.Item(71, 2).Value = 2726
$ws.Cells.Item(71, 3).Value = 10
$ws.Cells.Item(71, 4).Value = 1374
$ws.Cells.Item(71, 5).Value = 1201
$ws.Cells.Item(71, 6).Value = 32
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = 151

# Row 97: Consejo Danes para los Refugiados (was Mayotte; now ranks here)
$ws.Cells.Item(97, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(97, 2).Value = 1024
$ws.Cells.Item(97, 3).Value = 33
$ws.Cells.Item(97, 4).Value = 141
$ws.Cells.Item(97, 5).Value = 842
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 41

# Row 98: Mayotte (was Kirguistan; keeps the old Mayotte figures)
$ws.Cells.Item(98, 1).Value = "Mayotte"
$ws.Cells.Item(98, 2).Value = 1023
$ws.Cells.Item(98, 3).Value = 0
$ws.Cells.Item(98, 4).Value = 492
$ws.Cells.Item(98, 5).Value = 520
$ws.Cells.Item(98, 6).Value = 9
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 11

# Row 99: Kirguistan (was Consejo Danes; keeps the old Kirguistan figures)
$ws.Cells.Item(99, 1).Value = "Kirguistan"
$ws.Cells.Item(99, 2).Value = 1016
$ws.Cells.Item(99, 3).Value = 14
$ws.Cells.Item(99, 4).Value = 688
$ws.Cells.Item(99, 5).Value = 316
$ws.Cells.Item(99, 6).Value = 13
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 12

# Row 102: Republica de Chipre
$ws.Cells.Item(102, 1).Value = "Republica de Chipre"
$ws.Cells.Item(102, 2).Value = 901
$ws.Cells.Item(102, 3).Value = 3
$ws.Cells.Item(102, 4).Value = 401
$ws.Cells.Item(102, 5).Value = 484
$ws.Cells.Item(102, 6).Value = 10
$ws.Cells.Item(102, 7).Value = 0
$ws.Cells.Item(102, 8).Value = 16

# Row 106: Libano
$ws.Cells.Item(106, 1).Value = "Libano"
$ws.Cells.Item(106, 2).Value = 859
$ws.Cells.Item(106, 3).Value = 14
$ws.Cells.Item(106, 4).Value = 234
$ws.Cells.Item(106, 5).Value = 599
$ws.Cells.Item(106, 6).Value = 3
$ws.Cells.Item(106, 7).Value = 0
$ws.Cells.Item(106, 8).Value = 26

# Row 218: San Pedro y Miquelon
$ws.Cells.Item(218, 1).Value = "San Pedro y Miquelon"
$ws.Cells.Item(218, 2).Value = 1
$ws.Cells.Item(218, 3).Value = 0
$ws.Cells.Item(218, 4).Value = 1
$ws.Cells.Item(218, 5).Value = 0
$ws.Cells.Item(218, 6).Value = 0
$ws.Cells.Item(218, 7).Value = 0
$ws.Cells.Item(218, 8).Value = 0
